$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data point gets inserted right after the two most-recent rows
# (rows 2-3), pushing the previously-existing rows 4-9 down to rows 5-10.
$ws.Rows.Item(4).Insert()

$ws.Cells.Item(4, 1).Value = 8
$ws.Cells.Item(4, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(4, 3).Value = "Coquimbo"
$ws.Cells.Item(4, 4).Value = 44425
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
$ws.Cells.Item(4, 5).Value = 4
$ws.Cells.Item(4, 6).Value = 100112026
$ws.Cells.Item(4, 7).Value = "Haba"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 400
$ws.Cells.Item(4, 11).Value = 11500
$ws.Cells.Item(4, 12).Value = 12000
$ws.Cells.Item(4, 13).Value = 11750
$ws.Cells.Item(4, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(4, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(4, 16).Value = 470
$ws.Cells.Item(4, 17).Value = 25
$ws.Cells.Item(4, 18).Value = "Hortaliza"
